$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.430.36"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "'3.421.36"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'586.24"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "'137.57"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("D7").Value = "'3.424.18"
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'7.25"
$ws.Range("E10").Value = "  -5.11%  "
$ws.Range("D11").Value = "'0.121"
$ws.Range("E11").Value = "  -8.82%  "
$ws.Range("D12").Value = "'0.375"
$ws.Range("E12").Value = "  -6.84%  "
$ws.Range("D13").Value = "'4.000.77"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("E14").Value = "  -9.50%  "
$ws.Range("D15").Value = "'26.24"
$ws.Range("E15").Value = "  -8.01%  "
$ws.Range("D16").Value = "'3.424.89"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "'65.376.72"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "'9.80"
$ws.Range("E19").Value = "  -10.06%  "
$ws.Range("D20").Value = "'5.86"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("D21").Value = "'13.62"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("D22").Value = "'391.06"
$ws.Range("E22").Value = "  -4.98%  "
$ws.Range("D23").Value = "'0.554"
$ws.Range("E23").Value = "  -6.29%  "
$ws.Range("D24").Value = "'73.26"
$ws.Range("E24").Value = "  -5.43%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'3.563.11"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'0.0000106"
$ws.Range("E27").Value = "  -6.95%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'7.15"
$ws.Range("E29").Value = "  -6.32%  "
$ws.Range("D30").Value = "'8.17"
$ws.Range("E30").Value = "  -8.76%  "
$ws.Range("D31").Value = "'2.22"
$ws.Range("E31").Value = "  -8.29%  "
$ws.Range("D32").Value = "'3.429.24"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "'0.144"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("D35").Value = "'22.90"
$ws.Range("E35").Value = "  -5.01%  "
$ws.Range("D36").Value = "'172.38"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "'6.84"
$ws.Range("E37").Value = "  -7.75%  "
$ws.Range("D38").Value = "'1.16"
$ws.Range("E38").Value = "  -6.20%  "
$ws.Range("D39").Value = "'1.46"
$ws.Range("E39").Value = "  -6.07%  "
$ws.Range("D40").Value = "'4.76"
$ws.Range("E40").Value = "  -7.78%  "
$ws.Range("D41").Value = "'0.0762"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("D42").Value = "'0.818"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "'43.60"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -11.78%  "
$ws.Range("D46").Value = "'1.60"
$ws.Range("E46").Value = "  -8.69%  "
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "  +3.35%  "
$ws.Range("D48").Value = "'22.40"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").Value = "'6.51"
$ws.Range("E49").Value = "  -7.39%  "
$ws.Range("D50").Value = "'2.08"
$ws.Range("E50").Value = "  -12.27%  "
$ws.Range("D51").Value = "'2.192.00"
$ws.Range("E51").Value = "  -6.52%  "
